$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 24, shifting existing rows 24-32 down to 25-33
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new price entry
$ws.Range("A24").Value = 10
$ws.Range("B24").Value = "Vega Modelo de Temuco"
$ws.Range("C24").Value = "La Araucanía"
$ws.Range("D24").Value = 44704
$ws.Range("D24").NumberFormat = $ws.Range("D25").NumberFormat
$ws.Range("E24").Value = 9
$ws.Range("F24").Value = 100112010
$ws.Range("G24").Value = "Achicoria"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 100
$ws.Range("K24").Value = 11000
$ws.Range("L24").Value = 11000
$ws.Range("M24").Value = 11000
$ws.Range("N24").Value = "$/caja 18 unidades"
$ws.Range("O24").Value = "Región Metropolitana"
$ws.Range("P24").Value = 611
$ws.Range("Q24").Value = 18
$ws.Range("R24").Value = "Hortaliza"
